$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: sale line item data
$ws.Range("A7").Value = 1

$ws.Range("C7").NumberFormat = "@"
$ws.Range("C7").Value = "COLOVATIL 30 F.C. TABS"

$ws.Range("H7").Value = "2:0"

$ws.Range("L7").Value = "1"

$ws.Range("N7").NumberFormat = "@"
$ws.Range("N7").Value = "63.00"

$ws.Range("P7").NumberFormat = "@"
$ws.Range("P7").Value = "63.0000"

$ws.Range("Q7").Value = "1:0"

# Row 8: totals
$ws.Range("P8").Value = 63

# Footer timestamp refresh
$ws.Range("A9").Value = "Wednesday, 27 August, 2025 9:33 AM"
